# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.113.03"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.456.70"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.68"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.25"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.478"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.82"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.408"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "4.050.77"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.62"
$ws.Range("E14").Value = "  -4.43%  "
$ws.Range("D15").Value = "3.458.29"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "63.134.44"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.46"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.42"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.40"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "3.590.30"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.993"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.09"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.10"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.39"
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.33"
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "31.88"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.02"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.93"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "3.493.87"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0769"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.91"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.38"
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "2.587.38"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.88"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.62"
$ws.Range("E50").Value = "  -4.94%  "
$ws.Range("E51").Value = "  +0.04%  "
